$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1603853333333333
$ws.Range("H2").Value = 0.481156
$ws.Range("I2").Value = 0.01032935781992836
$ws.Range("J2").Value = 0.01042870175281933
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1579376666666667
$ws.Range("N2").Value = 0.473813
$ws.Range("O2").Value = 0.05467876644486869
$ws.Range("P2").Value = 0.07340983674118848
$ws.Range("Q2").Value = 0.02533088531422222
$ws.Range("R2").Value = 0.227977967828
$ws.Range("S2").Value = 0.0005647965437613408
$ws.Range("T2").Value = 0.0007655692930970135

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1603853333333333
$ws.Range("H3").Value = 0.481156
$ws.Range("I3").Value = 0.01032935781992836
$ws.Range("J3").Value = 0.01042870175281933
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.498127
$ws.Range("N3").Value = 1.494381
$ws.Range("O3").Value = 0.1724539210166233
$ws.Range("P3").Value = 0.2315307204300726
$ws.Range("Q3").Value = 0.07989226493733334
$ws.Range("R3").Value = 0.7190303844360001
$ws.Range("S3").Value = 0.001781338257630365
$ws.Range("T3").Value = 0.002414564829980621

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1603853333333333
$ws.Range("H4").Value = 0.481156
$ws.Range("I4").Value = 0.01032935781992836
$ws.Range("J4").Value = 0.01042870175281933
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02135966666666667
$ws.Range("N4").Value = 0.064079
$ws.Range("O4").Value = 0.00739481752299059
$ws.Range("P4").Value = 0.0099280284174107
$ws.Range("Q4").Value = 0.003425777258222222
$ws.Range("R4").Value = 0.030831995324
$ws.Range("S4").Value = 0.00007638371620804612
$ws.Range("T4").Value = 0.0001035364473586911

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1603853333333333
$ws.Range("H5").Value = 0.481156
$ws.Range("I5").Value = 0.01032935781992836
$ws.Range("J5").Value = 0.01042870175281933
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.21104
$ws.Range("N5").Value = 4.422079999999999
$ws.Range("O5").Value = 0.7654724950155174
$ws.Range("P5").Value = 0.6851314144113283
$ws.Range("Q5").Value = 0.3546183874133332
$ws.Range("R5").Value = 2.12771032448
$ws.Range("S5").Value = 0.007906839302328607
$ws.Range("T5").Value = 0.007145031182383009

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.92301466666667
$ws.Range("H6").Value = 44.76904399999999
$ws.Range("I6").Value = 0.9610926076617912
$ws.Range("J6").Value = 0.970336039943066
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1579376666666667
$ws.Range("N6").Value = 0.473813
$ws.Range("O6").Value = 0.05467876644486869
$ws.Range("P6").Value = 0.07340983674118848
$ws.Range("Q6").Value = 2.356906116085778
$ws.Range("R6").Value = 21.212155044772
$ws.Range("S6").Value = 0.0525513582262289
$ws.Range("T6").Value = 0.07123221027631182

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.92301466666667
$ws.Range("H7").Value = 44.76904399999999
$ws.Range("I7").Value = 0.9610926076617912
$ws.Range("J7").Value = 0.970336039943066
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.498127
$ws.Range("N7").Value = 1.494381
$ws.Range("O7").Value = 0.1724539210166233
$ws.Range("P7").Value = 0.2315307204300726
$ws.Range("Q7").Value = 7.433556526862667
$ws.Range("R7").Value = 66.90200874176399
$ws.Range("S7").Value = 0.165744188651367
$ws.Range("T7").Value = 0.2246626023872818

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.92301466666667
$ws.Range("H8").Value = 44.76904399999999
$ws.Range("I8").Value = 0.9610926076617912
$ws.Range("J8").Value = 0.970336039943066
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02135966666666667
$ws.Range("N8").Value = 0.064079
$ws.Range("O8").Value = 0.00739481752299059
$ws.Range("P8").Value = 0.0099280284174107
$ws.Range("Q8").Value = 0.3187506189417778
$ws.Range("R8").Value = 2.868755570475999
$ws.Range("S8").Value = 0.007107104456354134
$ws.Range("T8").Value = 0.009633523778992523

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.92301466666667
$ws.Range("H9").Value = 44.76904399999999
$ws.Range("I9").Value = 0.9610926076617912
$ws.Range("J9").Value = 0.970336039943066
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.21104
$ws.Range("N9").Value = 4.422079999999999
$ws.Range("O9").Value = 0.7654724950155174
$ws.Range("P9").Value = 0.6851314144113283
$ws.Range("Q9").Value = 32.99538234858666
$ws.Range("R9").Value = 197.97229409152
$ws.Range("S9").Value = 0.7356899563278411
$ws.Range("T9").Value = 0.66480770350048

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.443735
$ws.Range("H10").Value = 0.88747
$ws.Range("I10").Value = 0.02857803451828042
$ws.Range("J10").Value = 0.01923525830411462
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1579376666666667
$ws.Range("N10").Value = 0.473813
$ws.Range("O10").Value = 0.05467876644486869
$ws.Range("P10").Value = 0.07340983674118848
$ws.Range("Q10").Value = 0.07008247051833333
$ws.Range("R10").Value = 0.42049482311
$ws.Range("S10").Value = 0.00156261167487845
$ws.Range("T10").Value = 0.001412057171779645

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Myoc"
$ws.Range("C11").Value = "Fzd3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.443735
$ws.Range("H11").Value = 0.88747
$ws.Range("I11").Value = 0.02857803451828042
$ws.Range("J11").Value = 0.01923525830411462
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.498127
$ws.Range("N11").Value = 1.494381
$ws.Range("O11").Value = 0.1724539210166233
$ws.Range("P11").Value = 0.2315307204300726
$ws.Range("Q11").Value = 0.221036384345
$ws.Range("R11").Value = 1.32621830607
$ws.Range("S11").Value = 0.004928394107625865
$ws.Range("T11").Value = 0.004453553212810195

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Myoc"
$ws.Range("C12").Value = "Fzd3"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.443735
$ws.Range("H12").Value = 0.88747
$ws.Range("I12").Value = 0.02857803451828042
$ws.Range("J12").Value = 0.01923525830411462
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02135966666666667
$ws.Range("N12").Value = 0.064079
$ws.Range("O12").Value = 0.00739481752299059
$ws.Range("P12").Value = 0.0099280284174107
$ws.Range("Q12").Value = 0.009478031688333333
$ws.Range("R12").Value = 0.05686819013
$ws.Range("S12").Value = 0.00021132935042841
$ws.Range("T12").Value = 0.0001909681910594851

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Myoc"
$ws.Range("C13").Value = "Fzd3"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.443735
$ws.Range("H13").Value = 0.88747
$ws.Range("I13").Value = 0.02857803451828042
$ws.Range("J13").Value = 0.01923525830411462
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.21104
$ws.Range("N13").Value = 4.422079999999999
$ws.Range("O13").Value = 0.7654724950155174
$ws.Range("P13").Value = 0.6851314144113283
$ws.Range("Q13").Value = 0.9811158343999998
$ws.Range("R13").Value = 3.924463337599999
$ws.Range("S13").Value = 0.02187569938534769
$ws.Range("T13").Value = 0.0131786797284653
